$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title heading
#    ("Play Dragon & Phoenix Slot Game for Free - Betsoft"): a bold
#    "Meta description" label followed by the (non-bold) description text.
#
#    To make sure the new paragraph ends up with the same kind of paragraph
#    structure used throughout the rest of the document (Normal style, with
#    a leading empty run) we duplicate an existing Normal paragraph via
#    copy/paste and then swap its text for the text we actually want.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$normalSourcePara = $d.Paragraphs.Item(3)
$normalSourceText = $normalSourcePara.Range.Text

$normalSourcePara.Range.Copy()
$titlePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Paste()
$newPara.Style = $d.Styles.Item("Normal")

$boldLabel = "Meta description"
$restOfLine = ": Read our review of Dragon & Phoenix, the Chinese-themed online slot game from Betsoft. Play for free and enjoy this high-risk, high-reward game."

# Locate the pasted placeholder text so we know exactly where it starts/ends.
$placeholderRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $normalSourceText.Length)
$placeholderRange.Find.ClearFormatting()
$placeholderRange.Find.Execute($normalSourceText) | Out-Null

$oldStart = $placeholderRange.Start
$oldEnd = $placeholderRange.End

# Insert the new text right before the paragraph mark (oldEnd - 1 is always a
# safe, in-paragraph position, unlike oldEnd which sits at the start of the
# following paragraph).
$insPoint = $d.Range($oldEnd - 1, $oldEnd - 1)
$insPoint.InsertAfter($boldLabel + $restOfLine)

$boldRange = $d.Range($oldEnd - 1, $oldEnd - 1 + $boldLabel.Length)
$boldRange.Bold = 1

$restRange = $d.Range($oldEnd - 1 + $boldLabel.Length, $oldEnd - 1 + $boldLabel.Length + $restOfLine.Length)
$restRange.Bold = 0

# Remove the old placeholder text (but not the paragraph mark).
$delRange = $d.Range($oldStart, $oldEnd - 1)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Dragon & Phoenix Slot Game for Free -
#    Betsoft" paragraph that used to sit right before the italic summary
#    paragraph near the end of the document. That paragraph is always the
#    second-to-last paragraph in the document at this point.
# ---------------------------------------------------------------------------

$countAfterStep1 = $d.Paragraphs.Count
$dupPara = $d.Paragraphs.Item($countAfterStep1 - 1)
$dupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final italic paragraph with the new DALLE
#    image-prompt copy, while keeping its italic formatting.
# ---------------------------------------------------------------------------

$countAfterStep2 = $d.Paragraphs.Count
$summaryPara = $d.Paragraphs.Item($countAfterStep2)
$sumOldStart = $summaryPara.Range.Start
$sumOldEnd = $summaryPara.Range.End

$newSummary = 'DALLE, please create a feature image for the game "Dragon and Phoenix" that fits the theme and style of the game. The image should be in a cartoon style and should feature a happy Maya warrior with glasses. Make sure the image incorporates elements of the game such as the dragon, the Phoenix bird, the golden money tree, and the Emperor and Empress symbols. The image should be eye-catching and vibrant, and it should make people excited to play the game.'

$insPoint3 = $d.Range($sumOldEnd - 1, $sumOldEnd - 1)
$insPoint3.InsertAfter($newSummary)

$newSummaryRange = $d.Range($sumOldEnd - 1, $sumOldEnd - 1 + $newSummary.Length)
$newSummaryRange.Italic = 1

$oldSummaryRange = $d.Range($sumOldStart, $sumOldEnd - 1)
$oldSummaryRange.Delete()
